# "Data added for demo"
# Update the AUTOMATION customer-code demo data on the "addCustomer" sheet,
# adjust the active selection and widen the first column so the new values
# keep fitting, and nudge the workbook window size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addCustomer")

# --- Update the demo data values (row 2: CustomerName / Customer Code / Customer Pentacode) ---
$ws.Range("A2").Value = "AUTOMATION555"
$ws.Range("B2").Value = "AUTOMATION166"
$ws.Range("C2").Value = "AUTOMATION11555"

# --- Move/selection the active cell to A2 ---
$ws.Range("A2").Select()

# --- Widen column A so the longer values still fit ---
$ws.Columns.Item(1).ColumnWidth = 15

# --- Resize the workbook window ---
$excel.ActiveWindow.Width = 14280
$excel.ActiveWindow.Height = 6960
